# Insert two new rows at row 279 (shifts existing rows 279-360 down to 281-362)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("279:280").Insert()

# Populate the two newly inserted rows with the new daily quotations (date 44463 = 2021-09-24)
# Row 279
$ws.Cells.Item(279, 1).Value = 5
$ws.Cells.Item(279, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(279, 3).Value = 'Maule'
$ws.Cells.Item(279, 4).Value = 44463
$ws.Cells.Item(279, 5).Value = 7
$ws.Cells.Item(279, 6).Value = 'Fruta'
$ws.Cells.Item(279, 7).Value = 100108
$ws.Cells.Item(279, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(279, 9).Value = 100108006
$ws.Cells.Item(279, 10).Value = 'Plátano'
$ws.Cells.Item(279, 11).Value = 'Sin especificar'
$ws.Cells.Item(279, 12).Value = 'Pintón'
$ws.Cells.Item(279, 13).Value = 700
$ws.Cells.Item(279, 14).Value = 13000
$ws.Cells.Item(279, 15).Value = 13000
$ws.Cells.Item(279, 16).Value = 13000
$ws.Cells.Item(279, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(279, 18).Value = 'Ecuador'
$ws.Cells.Item(279, 19).Value = 650
$ws.Cells.Item(279, 20).Value = 20

# Row 280
$ws.Cells.Item(280, 1).Value = 5
$ws.Cells.Item(280, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(280, 3).Value = 'Maule'
$ws.Cells.Item(280, 4).Value = 44463
$ws.Cells.Item(280, 5).Value = 7
$ws.Cells.Item(280, 6).Value = 'Fruta'
$ws.Cells.Item(280, 7).Value = 100108
$ws.Cells.Item(280, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(280, 9).Value = 100108006
$ws.Cells.Item(280, 10).Value = 'Plátano'
$ws.Cells.Item(280, 11).Value = 'Sin especificar'
$ws.Cells.Item(280, 12).Value = 'Primera Pintón'
$ws.Cells.Item(280, 13).Value = 800
$ws.Cells.Item(280, 14).Value = 14000
$ws.Cells.Item(280, 15).Value = 15000
$ws.Cells.Item(280, 16).Value = 14375
$ws.Cells.Item(280, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(280, 18).Value = 'Ecuador'
$ws.Cells.Item(280, 19).Value = 719
$ws.Cells.Item(280, 20).Value = 20

